$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Összeadás")
$ws.Activate()

# Update F3 value and clear G3 entirely
$ws.Range("F3").Value = "Pass"
$ws.Range("G3").ClearContents()

# Update the selection to F3
$ws.Range("F3").Select()
